# Actualización desde MV -datos-
# Adds six new daily rows (09-10-2021 .. 14-10-2021) to the "Programa especial
# de compra de activos en pesos 2021 - Diaria" sheet, and backfills the B/C
# values (146 / 494) that were still pending for the 08-10-2021 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill row 282 (08-10-2021): it only had A/D/E filled in; add B/C ---
$ws.Range("B282").Value = 146
$ws.Range("C282").Value = 494

# --- New rows 283-287: full data (Serie/BCP/BCU/Bonos UF/Bonos $) ---
$newDates = @("09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")

$firstNewRow = 283
$lastFullRow = 287   # row 288 (14-10-2021) keeps the B/C columns empty, like 282 used to

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $firstNewRow + $i
    $dateText = $newDates[$i]

    # Write the literal date text without letting Excel auto-convert a
    # dd-mm-yyyy-looking string into a real date serial: build it as a text
    # formula, then paste-special as values so the cell collapses down to a
    # plain shared-string cell (same as every other "Serie" cell).
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = "=""" + $dateText + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)

    if ($row -le $lastFullRow) {
        $ws.Cells.Item($row, 2).Value = 146
        $ws.Cells.Item($row, 3).Value = 494
    }

    $ws.Cells.Item($row, 4).Value = 3088
    $ws.Cells.Item($row, 5).Value = 24
}

$excel.Application.CutCopyMode = $false
